$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Part 1: insert a new "Meta description" paragraph right after the
# top "Play Big Fin Bay Slot for Free - Exciting Marlin Hunting Theme"
# (Heading1) paragraph.
# ---------------------------------------------------------------------

$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphAfter()
$newp = $d.Paragraphs.Item(2)
$newp.Style = "Normal"

# The paragraph near the end of the document that currently reads
# "Play Big Fin Bay Slot for Free - Exciting Marlin Hunting Theme" in
# bold has the exact run layout we want to reproduce here (a leading
# empty run followed by a bold run). Clone its formatted content so the
# new paragraph picks up the same run structure, then retarget the text.
$count = $d.Paragraphs.Count
$srcPara = $d.Paragraphs.Item($count - 1)
$srcFT = $srcPara.Range.FormattedText
$newp.Range.FormattedText = $srcFT

$newp2 = $d.Paragraphs.Item(2)
$boldTextRange = $d.Range($newp2.Range.Start, $newp2.Range.End - 1)
$boldTextRange.Text = "Meta description"

$newp3 = $d.Paragraphs.Item(2)
$insertPos = $d.Range($newp3.Range.End - 1, $newp3.Range.End - 1)
$insertPos.InsertAfter(": Read our review of Big Fin Bay Slot and play for free with an expanding Wild symbol, free spins, and 117,649 ways to win during bonus rounds.")

# ---------------------------------------------------------------------
# Part 2: the duplicated bold "Play Big Fin Bay Slot for Free..."
# paragraph near the end of the document is removed, and the final
# italic paragraph's text is replaced with the new image prompt (its
# italic formatting is preserved since we only touch the text).
# ---------------------------------------------------------------------

$count2 = $d.Paragraphs.Count
$boldDup = $d.Paragraphs.Item($count2 - 1)
$boldDup.Range.Delete()

$count3 = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count3)
$lastTextRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$lastTextRange.Text = 'Prompt: Create a feature image for "Big Fin Bay" that showcases the adventure of the old sea wolf and the Marlin, while still reflecting the overall cartoon-style theme of the game. The image should feature the game''s main character, a happy Maya warrior with glasses, as he sets sail on his fishing boat with an optimistic grin on his face. The background should depict the beautiful, enchanting bay, with glimpses of the Marlin jumping out of the water. The image should be bright and colorful, with a playful tone that captures the excitement and thrill of the game.'
